# IONQ.xlsx update -- "JNJ and other updates"
#
# Main sheet:
#   - Price (J2) updated 7 -> 31.58 (J4/J7 formulas recalc automatically)
#   - New row 16: I16 = 100*2^14
#   - View: zoom 250 -> 175, selection moved to K5
#
# Model sheet:
#   - New quarterly datapoints across row 3 (F3:K3) plus M3, alongside the
#     pre-existing L3 value
#   - View: bottom-right pane selection moved from L6 to M3

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Main sheet
# ---------------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("Main")

$wsMain.Range("J2").Value = 31.58

$wsMain.Range("I16").Formula = "=100*2^14"

# ---------------------------------------------------------------------------
# Model sheet
# ---------------------------------------------------------------------------
$wsModel = $wb.Worksheets.Item("Model")

$wsModel.Range("F3").Value = 3.8
$wsModel.Range("G3").Value = 4.3
$wsModel.Range("H3").Value = 5.5
$wsModel.Range("I3").Value = 6.1
$wsModel.Range("J3").Value = 6.1
$wsModel.Range("K3").Value = 7.6
$wsModel.Range("M3").Value = 12.4

# Update Model's pane selection first ...
$wsModel.Activate()
$wsModel.Range("M3").Select()

# ... then re-activate Main (keeps Main as the selected/visible tab, matches
# the original file) and update its zoom + selection.
$wsMain.Activate()
$excel.ActiveWindow.Zoom = 175
$wsMain.Range("K5").Select()
